$d = $word.ActiveDocument

$replacements = @(
    @("484÷8=", "358÷3="),
    @("344÷7=", "238÷7="),
    @("974÷7=", "646÷2="),
    @("549÷4=", "715÷6="),
    @("179÷9=", "298÷3="),
    @("753÷5=", "337÷3="),
    @("306÷6=", "807÷2="),
    @("665÷8=", "730÷8="),
    @("890÷5=", "756÷3="),
    @("753÷8=", "961÷6="),
    @("678÷7=", "565÷2="),
    @("585÷6=", "346÷4="),
    @("462÷4=", "820÷4="),
    @("784÷4=", "104÷6="),
    @("834÷3=", "964÷9="),
    @("135÷8=", "566÷8="),
    @("558÷8=", "997÷8="),
    @("684÷2=", "390÷2="),
    @("842÷8=", "932÷7="),
    @("974÷4=", "937÷6="),
    @("328÷8=", "656÷2="),
    @("291÷9=", "173÷9="),
    @("787÷2=", "299÷6="),
    @("988÷4=", "209÷4="),
    @("428÷8=", "429÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
